$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4,D5,D6,D7,D9,D10,D11,D14,D15,D17,D21,D22,D23,D24,D26,D27,D28,D29,D30,D31,D33,D35,D37,D38,D39,D40,D41,D44,D45,D46,D49,D50,D51").NumberFormat = "@"

$ws.Range("D2").Value = "40.043.99"
$ws.Range("E2").Value = "  -4.02%  "
$ws.Range("D3").Value = "2.330.90"
$ws.Range("E3").Value = "  -5.67%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "308.51"
$ws.Range("E5").Value = "  -3.90%  "
$ws.Range("D6").Value = "84.94"
$ws.Range("E6").Value = "  -7.89%  "
$ws.Range("D7").Value = "0.532"
$ws.Range("E7").Value = "  -3.18%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "0.485"
$ws.Range("E9").Value = "  -4.46%  "
$ws.Range("D10").Value = "0.0814"
$ws.Range("E10").Value = "  -4.68%  "
$ws.Range("D11").Value = "30.06"
$ws.Range("E11").Value = "  -8.66%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "2.691.22"
$ws.Range("E13").Value = "  -5.68%  "
$ws.Range("D14").Value = "6.41"
$ws.Range("E14").Value = "  -6.73%  "
$ws.Range("D15").Value = "14.71"
$ws.Range("E15").Value = "  -5.09%  "
$ws.Range("D16").Value = "2.331.41"
$ws.Range("E16").Value = "  -5.64%  "
$ws.Range("D17").Value = "0.755"
$ws.Range("E17").Value = "  -4.30%  "
$ws.Range("D18").Value = "39.999.54"
$ws.Range("E18").Value = "  -4.02%  "
$ws.Range("D19").Value = "0.0₃0903"
$ws.Range("E19").Value = "  -4.03%  "
$ws.Range("E20").Value = "  -5.69%  "
$ws.Range("D21").Value = "67.63"
$ws.Range("E21").Value = "  -6.00%  "
$ws.Range("D22").Value = "10.63"
$ws.Range("E22").Value = "  -5.16%  "
$ws.Range("D23").Value = "235.72"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").Value = "2.56"
$ws.Range("E24").Value = "  -7.05%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "1.80"
$ws.Range("E26").Value = "  -6.70%  "
$ws.Range("D27").Value = "23.35"
$ws.Range("E27").Value = "  -5.78%  "
$ws.Range("D28").Value = "2.14"
$ws.Range("E28").Value = "  -4.38%  "
$ws.Range("D29").Value = "9.27"
$ws.Range("E29").Value = "  -4.43%  "
$ws.Range("D30").Value = "35.21"
$ws.Range("E30").Value = "  -2.59%  "
$ws.Range("D31").Value = "152.49"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "5.11"
$ws.Range("E33").Value = "  -5.81%  "
$ws.Range("E34").Value = "  -4.22%  "
$ws.Range("D35").Value = "0.0722"
$ws.Range("E35").Value = "  -5.23%  "
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("D37").Value = "0.0999"
$ws.Range("E37").Value = "  -3.07%  "
$ws.Range("D38").Value = "2.75"
$ws.Range("E38").Value = "  -4.92%  "
$ws.Range("D39").Value = "15.66"
$ws.Range("E39").Value = "  -7.95%  "
$ws.Range("D40").Value = "1.71"
$ws.Range("E40").Value = "  -6.81%  "
$ws.Range("D41").Value = "3.83"
$ws.Range("E41").Value = "  -3.77%  "
$ws.Range("E42").Value = "  -3.41%  "
$ws.Range("D43").Value = "1.951.28"
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("D44").Value = "0.0266"
$ws.Range("E44").Value = "  -5.48%  "
$ws.Range("D45").Value = "17.55"
$ws.Range("E45").Value = "  -5.80%  "
$ws.Range("D46").Value = "9.35"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("E47").Value = "  -9.34%  "
$ws.Range("D48").Value = "2.552.98"
$ws.Range("E48").Value = "  -6.63%  "
$ws.Range("D49").Value = "92.80"
$ws.Range("E49").Value = "  -4.58%  "
$ws.Range("D50").Value = "70.52"
$ws.Range("E50").Value = "  -6.94%  "
$ws.Range("D51").Value = "50.42"
$ws.Range("E51").Value = "  -3.25%  "
